# Add new CRM rows for three new contacts (Kyro Johnsen, Ryderi Reamer,
# Patrick-John Caswell) to Sheet1 (deal/contact info) and Sheet2
# (phone-number detail rows), matching the upstream commit "Add files via
# upload".
#
# NOTE: cell writes are interleaved in the exact sequence in which the
# corresponding shared-string values were first introduced upstream, so
# that the rebuilt sharedStrings table lines up with the target file.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

function Set-DateCell($cell, $y, $m, $d) {
    $cell.Value = (Get-Date -Year $y -Month $m -Day $d -Hour 0 -Minute 0 -Second 0)
}

function Set-TextNumberCell($cell, $n) {
    # Sheet2 column C is styled with a Text ("@") number format, yet the
    # source workbook stores a genuine numeric value (no shared-string
    # reference) in it. Temporarily switch the cell to the default
    # "Normal" style (which maps to the existing General cellXf) so the
    # numeric literal is written as a number, then restore the original
    # Text format without allocating any new style entries.
    $cell.Style = "Normal"
    $cell.Value = $n
    $cell.NumberFormat = "@"
}

# ---------------------------------------------------------------------
# Person 1: Kyro Johnsen (row 118 in Sheet1; deals C4J84051MW /
# O78YEEH16E / X24EU69882, each with two Sheet2 phone rows)
# ---------------------------------------------------------------------

$ws1.Cells.Item(118,1).Value = "342000655"
$ws1.Cells.Item(118,2).Value = "Kyro"
$ws1.Cells.Item(118,3).Value = "Johnsen"
$ws1.Cells.Item(118,4).Value = "0527290650"
$ws1.Cells.Item(118,5).Value = "Kyro_Johnsen@Jmail.com"
$ws1.Cells.Item(118,6).Value = "Closed lost(0%)"
$ws1.Cells.Item(118,7).Value = "204264543"
Set-DateCell $ws1.Cells.Item(118,8) 2018 12 11
$ws1.Cells.Item(118,9).Value = "don" + [char]0x2019 + "t know what he wants in his life"
$ws1.Cells.Item(118,10).Value = "C4J84051MW"

$ws2.Cells.Item(259,1).Value = "342000655"
$ws2.Cells.Item(259,2).Value = "0553433342"
Set-TextNumberCell $ws2.Cells.Item(259,3) 3
Set-DateCell $ws2.Cells.Item(259,4) 2018 11 12
Set-DateCell $ws2.Cells.Item(259,5) 2019 5 12
$ws2.Cells.Item(259,6).Value = "C4J84051MW"

$ws2.Cells.Item(260,1).Value = "342000655"
$ws2.Cells.Item(260,2).Value = "0536921410"
Set-TextNumberCell $ws2.Cells.Item(260,3) 2
Set-DateCell $ws2.Cells.Item(260,4) 2018 11 13
Set-DateCell $ws2.Cells.Item(260,5) 2019 7 13
$ws2.Cells.Item(260,6).Value = "C4J84051MW"

$ws1.Cells.Item(119,1).Value = "342000655"
$ws1.Cells.Item(119,2).Value = "Kyro"
$ws1.Cells.Item(119,3).Value = "Johnsen"
$ws1.Cells.Item(119,4).Value = "0527290650"
$ws1.Cells.Item(119,5).Value = "Kyro_Johnsen@Jmail.com"
$ws1.Cells.Item(119,6).Value = "Closed lost(0%)"
$ws1.Cells.Item(119,7).Value = "204264543"
Set-DateCell $ws1.Cells.Item(119,8) 2018 12 11
$ws1.Cells.Item(119,9).Value = "don" + [char]0x2019 + "t know what he wants in his life"
$ws1.Cells.Item(119,10).Value = "O78YEEH16E"

$ws2.Cells.Item(261,1).Value = "342000655"
$ws2.Cells.Item(261,2).Value = "0553433342"
Set-TextNumberCell $ws2.Cells.Item(261,3) 3
Set-DateCell $ws2.Cells.Item(261,4) 2018 11 12
Set-DateCell $ws2.Cells.Item(261,5) 2019 5 12
$ws2.Cells.Item(261,6).Value = "O78YEEH16E"

$ws2.Cells.Item(262,1).Value = "342000655"
$ws2.Cells.Item(262,2).Value = "0536921410"
Set-TextNumberCell $ws2.Cells.Item(262,3) 2
Set-DateCell $ws2.Cells.Item(262,4) 2018 11 13
Set-DateCell $ws2.Cells.Item(262,5) 2019 7 13
$ws2.Cells.Item(262,6).Value = "O78YEEH16E"

$ws1.Cells.Item(120,1).Value = "342000655"
$ws1.Cells.Item(120,2).Value = "Kyro"
$ws1.Cells.Item(120,3).Value = "Johnsen"
$ws1.Cells.Item(120,4).Value = "0527290650"
$ws1.Cells.Item(120,5).Value = "Kyro_Johnsen@Jmail.com"
$ws1.Cells.Item(120,6).Value = "Closed lost(0%)"
$ws1.Cells.Item(120,7).Value = "204264543"
Set-DateCell $ws1.Cells.Item(120,8) 2018 12 11
$ws1.Cells.Item(120,9).Value = "don" + [char]0x2019 + "t know what he wants in his life"
$ws1.Cells.Item(120,10).Value = "X24EU69882"

$ws2.Cells.Item(263,1).Value = "342000655"
$ws2.Cells.Item(263,2).Value = "0553433342"
Set-TextNumberCell $ws2.Cells.Item(263,3) 3
Set-DateCell $ws2.Cells.Item(263,4) 2018 11 12
Set-DateCell $ws2.Cells.Item(263,5) 2019 5 12
$ws2.Cells.Item(263,6).Value = "X24EU69882"

$ws2.Cells.Item(264,1).Value = "342000655"
$ws2.Cells.Item(264,2).Value = "0536921410"
Set-TextNumberCell $ws2.Cells.Item(264,3) 2
Set-DateCell $ws2.Cells.Item(264,4) 2018 11 13
Set-DateCell $ws2.Cells.Item(264,5) 2019 7 13
$ws2.Cells.Item(264,6).Value = "X24EU69882"

# ---------------------------------------------------------------------
# Person 2: Ryderi Reamer (row 121 in Sheet1; single deal 4PQCJZGZ7Q
# with three Sheet2 phone rows)
# ---------------------------------------------------------------------

$ws1.Cells.Item(121,1).Value = "171073405"
$ws1.Cells.Item(121,2).Value = "Ryderi"
$ws1.Cells.Item(121,3).Value = "Reamer"
$ws1.Cells.Item(121,4).Value = "0583113772"
$ws1.Cells.Item(121,5).Value = "Ryder_Reamer@Lmail.com"
$ws1.Cells.Item(121,6).Value = "Closed lost(0%)"
$ws1.Cells.Item(121,7).Value = "208063511"
Set-DateCell $ws1.Cells.Item(121,8) 2018 12 11
$ws1.Cells.Item(121,9).Value = "don" + [char]0x2019 + "t know what he wants in his life yet"
$ws1.Cells.Item(121,10).Value = "4PQCJZGZ7Q"

$ws2.Cells.Item(265,1).Value = "171073405"
$ws2.Cells.Item(265,2).Value = "0533744671"
Set-TextNumberCell $ws2.Cells.Item(265,3) 1
Set-DateCell $ws2.Cells.Item(265,4) 2018 11 12
Set-DateCell $ws2.Cells.Item(265,5) 2019 6 12
$ws2.Cells.Item(265,6).Value = "4PQCJZGZ7Q"

$ws2.Cells.Item(266,1).Value = "171073405"
$ws2.Cells.Item(266,2).Value = "0581279283"
Set-TextNumberCell $ws2.Cells.Item(266,3) 2
Set-DateCell $ws2.Cells.Item(266,4) 2018 11 9
Set-DateCell $ws2.Cells.Item(266,5) 2019 11 9
$ws2.Cells.Item(266,6).Value = "4PQCJZGZ7Q"

$ws2.Cells.Item(267,1).Value = "171073405"
$ws2.Cells.Item(267,2).Value = "0506954641"
Set-TextNumberCell $ws2.Cells.Item(267,3) 3
Set-DateCell $ws2.Cells.Item(267,4) 2018 11 13
Set-DateCell $ws2.Cells.Item(267,5) 2019 8 13
$ws2.Cells.Item(267,6).Value = "4PQCJZGZ7Q"

# ---------------------------------------------------------------------
# Person 3: Patrick-John Caswell (row 122 in Sheet1; single deal
# 94OU5KH91Q with two Sheet2 phone rows)
# ---------------------------------------------------------------------

$ws1.Cells.Item(122,1).Value = "236263121"
$ws1.Cells.Item(122,2).Value = "Patrick-John"
$ws1.Cells.Item(122,3).Value = "Caswell"
$ws1.Cells.Item(122,4).Value = "0526700677"
$ws1.Cells.Item(122,5).Value = "Patrick-John_Caswell@Umail.com"
$ws1.Cells.Item(122,6).Value = "Closed lost(0%)"
$ws1.Cells.Item(122,7).Value = "204264543"
Set-DateCell $ws1.Cells.Item(122,8) 2018 12 11
$ws1.Cells.Item(122,9).Value = "don" + [char]0x2019 + "t know what he wants in his life"
$ws1.Cells.Item(122,10).Value = "94OU5KH91Q"

$ws2.Cells.Item(268,1).Value = "236263121"
$ws2.Cells.Item(268,2).Value = "0517142261"
Set-TextNumberCell $ws2.Cells.Item(268,3) 3
Set-DateCell $ws2.Cells.Item(268,4) 2018 11 4
Set-DateCell $ws2.Cells.Item(268,5) 2019 6 4
$ws2.Cells.Item(268,6).Value = "94OU5KH91Q"

$ws2.Cells.Item(269,1).Value = "236263121"
$ws2.Cells.Item(269,2).Value = "0523344409"
Set-TextNumberCell $ws2.Cells.Item(269,3) 2
Set-DateCell $ws2.Cells.Item(269,4) 2018 11 7
Set-DateCell $ws2.Cells.Item(269,5) 2019 11 7
$ws2.Cells.Item(269,6).Value = "94OU5KH91Q"
